# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de):
#   - Status column (B) flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" for the two tracked files.
#   - Two new columns get populated for those same two rows:
#       E "Latest Target File"   -> same file as the source markdown (A)
#       F "Latest Handback File" -> same file as the latest handoff xlf (C)
#     both rendered as hyperlinks mirroring A/C's links.
#   - Column G "Latest Handback DateTime" is stamped with the handback time.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$sheetInfo = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-03-08 02:44:48" },
    @{ Name = "de-de"; HandbackTime = "2016-03-08 02:45:02" }
)

foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    foreach ($row in @(2, 3)) {
        $aCell = $ws.Cells.Item($row, 1)   # A: Source File Name
        $cCell = $ws.Cells.Item($row, 3)   # C: Latest Handoff File
        $eCell = $ws.Cells.Item($row, 5)   # E: Latest Target File
        $fCell = $ws.Cells.Item($row, 6)   # F: Latest Handback File
        $gCell = $ws.Cells.Item($row, 7)   # G: Latest Handback DateTime

        # Status -> handed back
        $ws.Cells.Item($row, 2).Value2 = $newStatus

        # Find the existing hyperlink addresses for A and C on this row so the
        # new Target File / Handback File links point at the same place.
        $aAddress = $null
        $cAddress = $null
        foreach ($hl in $ws.Hyperlinks) {
            if ($hl.Range.Address() -eq $aCell.Address()) { $aAddress = $hl.Address }
            if ($hl.Range.Address() -eq $cCell.Address()) { $cAddress = $hl.Address }
        }

        $aText = $aCell.Value2
        $cText = $cCell.Value2

        $ws.Hyperlinks.Add($eCell, $aAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $aText) | Out-Null
        $ws.Hyperlinks.Add($fCell, $cAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $cText) | Out-Null

        # Handback timestamp
        $gCell.Value2 = $info.HandbackTime
    }
}
